$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextCell $ws.Range("D2") "26.973.51"
Set-TextCell $ws.Range("E2") "  -3.05%  "
Set-TextCell $ws.Range("D3") "1.737.43"
Set-TextCell $ws.Range("E3") "  -1.30%  "
Set-TextCell $ws.Range("D4") "1.001"
Set-TextCell $ws.Range("E4") "  -0.22%  "
Set-TextCell $ws.Range("D5") "310.97"
Set-TextCell $ws.Range("E5") "  -5.49%  "
Set-TextCell $ws.Range("D6") "0.9998"
Set-TextCell $ws.Range("E6") "  -0.27%  "
Set-TextCell $ws.Range("D7") "0.4997"
Set-TextCell $ws.Range("E7") "  +3.58%  "
Set-TextCell $ws.Range("E8") "  -0.84%  "
Set-TextCell $ws.Range("D9") "42.63"
Set-TextCell $ws.Range("D10") "0.07250"
Set-TextCell $ws.Range("E10") "  -3.67%  "
Set-TextCell $ws.Range("D11") "1.057"
Set-TextCell $ws.Range("E11") "  -1.97%  "
Set-TextCell $ws.Range("D12") "1.001"
Set-TextCell $ws.Range("E12") "  -0.17%  "
Set-TextCell $ws.Range("D13") "20.04"
Set-TextCell $ws.Range("E13") "  -2.55%  "
Set-TextCell $ws.Range("D14") "5.943"
Set-TextCell $ws.Range("E14") "  -1.74%  "
Set-TextCell $ws.Range("D15") "1.731.78"
Set-TextCell $ws.Range("E15") "  -1.76%  "
Set-TextCell $ws.Range("D16") "6.875"
Set-TextCell $ws.Range("E16") "  -3.97%  "
Set-TextCell $ws.Range("D17") "86.34"
Set-TextCell $ws.Range("E17") "  -6.52%  "
Set-TextCell $ws.Range("E18") "  -4.48%  "
Set-TextCell $ws.Range("D19") "0.06385"
Set-TextCell $ws.Range("E19") "  -0.89%  "
Set-TextCell $ws.Range("D20") "0.9994"
Set-TextCell $ws.Range("E20") "  -0.21%  "
Set-TextCell $ws.Range("D21") "16.62"
Set-TextCell $ws.Range("E21") "  -1.34%  "
Set-TextCell $ws.Range("D22") "5.734"
Set-TextCell $ws.Range("E22") "  -0.64%  "
Set-TextCell $ws.Range("D23") "27.068.28"
Set-TextCell $ws.Range("E23") "  -2.84%  "
Set-TextCell $ws.Range("D24") "11.05"
Set-TextCell $ws.Range("E24") "  -0.60%  "
Set-TextCell $ws.Range("D25") "2.050"
Set-TextCell $ws.Range("E25") "  -5.13%  "
Set-TextCell $ws.Range("D26") "153.64"
Set-TextCell $ws.Range("E26") "  -6.30%  "
Set-TextCell $ws.Range("D27") "19.94"
Set-TextCell $ws.Range("E27") "  -0.62%  "
Set-TextCell $ws.Range("D28") "1.935.51"
Set-TextCell $ws.Range("E28") "  -1.49%  "
Set-TextCell $ws.Range("D29") "2.112"
Set-TextCell $ws.Range("E29") "  -3.33%  "
Set-TextCell $ws.Range("D30") "120.48"
Set-TextCell $ws.Range("E30") "  -2.22%  "
Set-TextCell $ws.Range("D31") "1.064"
Set-TextCell $ws.Range("E31") "  +0.14%  "
Set-TextCell $ws.Range("D32") "0.09499"
Set-TextCell $ws.Range("E32") "  +0.51%  "
Set-TextCell $ws.Range("D33") "3.576"
Set-TextCell $ws.Range("E33") "  -2.36%  "
Set-TextCell $ws.Range("D34") "5.388"
Set-TextCell $ws.Range("E34") "  -2.73%  "
Set-TextCell $ws.Range("D35") "0.05937"
Set-TextCell $ws.Range("E35") "  -1.86%  "
Set-TextCell $ws.Range("D36") "0.02191"
Set-TextCell $ws.Range("E36") "  -3.12%  "
Set-TextCell $ws.Range("D37") "11.03"
Set-TextCell $ws.Range("E37") "  -5.17%  "
Set-TextCell $ws.Range("D38") "1.433"
Set-TextCell $ws.Range("E38") "  -0.49%  "
Set-TextCell $ws.Range("B39") "Algorand"
Set-TextCell $ws.Range("C39") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell $ws.Range("D39") "0.1993"
Set-TextCell $ws.Range("E39") "  -3.14%  "
Set-TextCell $ws.Range("B40") "InternetComputer(DFINITY)"
Set-TextCell $ws.Range("C40") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell $ws.Range("D40") "4.761"
Set-TextCell $ws.Range("E40") "  -2.56%  "
Set-TextCell $ws.Range("B41") "Frax"
Set-TextCell $ws.Range("C41") "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextCell $ws.Range("D41") "0.9996"
Set-TextCell $ws.Range("E41") "  -0.26%  "
Set-TextCell $ws.Range("B42") "TheSandbox"
Set-TextCell $ws.Range("C42") "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextCell $ws.Range("D42") "0.6006"
Set-TextCell $ws.Range("E42") "  -2.20%  "
Set-TextCell $ws.Range("D43") "1.111"
Set-TextCell $ws.Range("E43") "  -5.82%  "
Set-TextCell $ws.Range("D44") "7.474"
Set-TextCell $ws.Range("E44") "  -3.23%  "
Set-TextCell $ws.Range("D45") "12.83"
Set-TextCell $ws.Range("E45") "  -2.29%  "
Set-TextCell $ws.Range("D46") "3.577"
Set-TextCell $ws.Range("E46") "  -4.35%  "
Set-TextCell $ws.Range("D47") "0.5638"
Set-TextCell $ws.Range("E47") "  -2.56%  "
Set-TextCell $ws.Range("D48") "119.45"
Set-TextCell $ws.Range("E48") "  -3.04%  "
Set-TextCell $ws.Range("D49") "1.852"
Set-TextCell $ws.Range("E49") "  -3.68%  "
Set-TextCell $ws.Range("E50") "  -2.04%  "
Set-TextCell $ws.Range("D51") "1.098"
Set-TextCell $ws.Range("E51") "  -3.56%  "
